$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells D2 through J2 currently hold "Unknown" and should become "unknown"
$ws.Range("D2:J2").Value = "unknown"
